# overview_testcases.xlsx — "removed xlsx cached stuff"
#
# Semantic changes applied:
#   1. overview_testcases!I39: 128 -> 512 (raw input; I40:I48 reference it via
#      "=$I$39" and K39:K48 recompute (MAX/AVERAGE ratio) automatically on
#      recalc).
#   2. Active sheet switches from "overhead" (2nd tab) back to
#      "overview_testcases" (1st tab).
#   3. overview_testcases sheet view: zoomed to 80%, selection moved to M43.
#   4. overhead sheet view: no longer the tab-selected sheet (selection on
#      G15 is left untouched).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("overview_testcases")
$wsOverhead = $wb.Worksheets.Item("overhead")

# 1. Data edit: bump the P6 column's raw sample for the "varied imbalance
# ratios" table from 128 to 512. Downstream formulas (=$I$39 copies and the
# MAX/AVERAGE "Rimb" ratios) recalc automatically.
$wsOverview.Range("I39").Value = 512

# 2 & 3. Make overview_testcases the active/visible tab again, zoom its view
# to 80%, and move the selection to M43.
$wsOverview.Activate()
$excel.ActiveWindow.Zoom = 80
$wsOverview.Range("M43").Select()
